# Auto-generated edit script for SCC Testcases sheet update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testcases")

$ws.Range("B2").Value = "Component: SCC (Service Cloud Client)"
$ws.Range("D6").Value = "Verify enabling/disabling of SCC function"
$ws.Range("E6").Value = "1. Enter Self-diagnostic boot-up mode`n2. Set 08-3820 to 1 (Enable)`n3. Reboot the MFP`n4. Check if SCC function is enabled"
$ws.Range("F6").Value = "SCC function should be enabled after reboot"
$ws.Range("H6").Value = "Default value should be `"Disable`""
$ws.Range("C7").Value = "MFP is powered on and in normal boot-up mode"
$ws.Range("D7").Value = "Verify SCC function cannot be enabled from administrator mode"
$ws.Range("E7").Value = "1. Login as administrator`n2. Try to access SCC settings from TopAccess`n3. Verify no SCC settings are available"
$ws.Range("F7").Value = "Administrator should not be able to enable/disable SCC function from TopAccess"
$ws.Range("C8").Value = "MFP with SCC disabled"
$ws.Range("D8").Value = "Verify SCC installation report printing"
$ws.Range("E8").Value = "1. Enter Self-diagnostic boot-up mode`n2. Change 08-3820 from 0 (Disable) to 1 (Enable)`n3. Reboot the MFP`n4. Wait for MFP registration processing to complete"
$ws.Range("F8").Value = "SCC installation report should be automatically printed after registration processing"
$ws.Range("H8").Value = "Report should print only once when SCC is enabled for the first time"
$ws.Range("C9").Value = "MFP with SCC enabled"
$ws.Range("D9").Value = "Verify SCC installation report content"
$ws.Range("E9").Value = "1. Check the printed installation report"
$ws.Range("F9").Value = "Report should contain: Serial Number, MAC Address, Model Name, Firmware Version, SCC Version, Registration Status, Registration Time, Server URL"
$ws.Range("H9").Value = "Report should be in English only"
$ws.Range("C10").Value = "MFP with SCC enabled and network print restriction mode enabled"
$ws.Range("D10").Value = "Verify SCC report printing with network print restrictions"
$ws.Range("E10").Value = "1. Set 08-9344 to 1 (only private)`n2. Enable SCC`n3. Check if report is printed"
$ws.Range("C11").Value = "MFP with SCC enabled and user authentication enabled"
$ws.Range("D11").Value = "Verify SCC report printing with quota restrictions"
$ws.Range("E11").Value = "1. Enable user authentication`n2. Set quota for built-in admin to 0`n3. Enable SCC`n4. Check if report is printed"
$ws.Range("F11").Value = "SCC report should not be printed due to quota error"
$ws.Range("C12").Value = "MFP with SCC enabled"
$ws.Range("D12").Value = "Verify proxy server settings"
$ws.Range("E12").Value = "1. Enter Self-diagnostic boot-up mode`n2. Configure proxy settings (08-3822 to 08-3826)`n3. Reboot MFP`n4. Verify SCC communication works through proxy"
$ws.Range("F12").Value = "SCC should communicate with server through configured proxy"
$ws.Range("C13").Value = "MFP with SCC enabled"
$ws.Range("D13").Value = "Verify communication protocol"
$ws.Range("E13").Value = "1. Monitor network traffic during SCC communication`n2. Verify HTTPS is used"
$ws.Range("F13").Value = "All SCC communication should use HTTPS protocol"
$ws.Range("C14").Value = "MFP with SCC enabled"
$ws.Range("D14").Value = "Verify port number fallback"
$ws.Range("E14").Value = "1. Block port 443 at firewall`n2. Trigger SCC communication`n3. Monitor network traffic"
$ws.Range("F14").Value = "SCC should automatically use port 8443 when port 443 is not available"
$ws.Range("C15").Value = "MFP with SCC enabled"
$ws.Range("D15").Value = "Verify SCC does not start in special startup mode"
$ws.Range("E15").Value = "1. Boot MFP in special startup mode`n2. Check if SCC process starts"
$ws.Range("F15").Value = "SCC process should not start in special startup mode"
$ws.Range("C16").Value = "MFP with SCC enabled"
$ws.Range("D16").Value = "Verify regular communication loop timing"
$ws.Range("E16").Value = "1. Set regular communication schedule`n2. Monitor when SCC communicates with server"
$ws.Range("F16").Value = "SCC should communicate with server according to configured schedule"
$ws.Range("H16").Value = "Default schedule is `"every day at 0:00`""
$ws.Range("C17").Value = "MFP with SCC enabled"
$ws.Range("D17").Value = "Verify first-time device registration"
$ws.Range("E17").Value = "1. Enable SCC for the first time`n2. Monitor registration process"
$ws.Range("F17").Value = "MFP should register with server and receive authentication token"
$ws.Range("C18").Value = "MFP with SCC registered"
$ws.Range("D18").Value = "Verify subsequent authentication"
$ws.Range("E18").Value = "1. Trigger regular communication after registration`n2. Monitor authentication process"
$ws.Range("F18").Value = "MFP should authenticate using token received during registration"
$ws.Range("C19").Value = "MFP with SCC enabled"
$ws.Range("D19").Value = "Verify server busy handling"
$ws.Range("E19").Value = "1. Simulate server BUSY response`n2. Monitor MFP behavior"
$ws.Range("F19").Value = "MFP should retry connection with increasing intervals as specified by server"
$ws.Range("H19").Value = "After 3 retries, MFP should exit communication cycle"
$ws.Range("C20").Value = "MFP with SCC enabled"
$ws.Range("D20").Value = "Verify registration failure handling"
$ws.Range("E20").Value = "1. Simulate registration failure response`n2. Monitor MFP behavior"
$ws.Range("F20").Value = "After 3 consecutive failures, MFP should print Installation Report and exit communication cycle"
$ws.Range("C21").Value = "MFP with SCC enabled"
$ws.Range("D21").Value = "Verify Check for Updates functionality"
$ws.Range("E21").Value = "1. Trigger regular communication`n2. Monitor Check for Updates request"
$ws.Range("F21").Value = "MFP should send correct parameters and process server response appropriately"
$ws.Range("C22").Value = "MFP with SCC enabled"
$ws.Range("D22").Value = "Verify Download Package functionality"
$ws.Range("E22").Value = "1. Configure server to send update package`n2. Trigger regular communication`n3. Monitor download process"
$ws.Range("F22").Value = "MFP should download package and verify hash value"
$ws.Range("C23").Value = "MFP with SCC enabled"
$ws.Range("D23").Value = "Verify package hash validation"
$ws.Range("E23").Value = "1. Configure server to send package with incorrect hash`n2. Trigger regular communication`n3. Monitor MFP behavior"
$ws.Range("F23").Value = "MFP should delete downloaded data and exit communication cycle"
$ws.Range("C24").Value = "MFP with SCC enabled"
$ws.Range("D24").Value = "Verify Send Baseline Data functionality"
$ws.Range("E24").Value = "1. Trigger regular communication`n2. Monitor baseline data transmission"
$ws.Range("F24").Value = "MFP should collect and send all required data according to SendDataConfig setting"
$ws.Range("C25").Value = "MFP with SCC enabled"
$ws.Range("D25").Value = "Verify Install Package for firmware update"
$ws.Range("E25").Value = "1. Configure server to send firmware update package`n2. Trigger regular communication`n3. Monitor installation process"
$ws.Range("F25").Value = "MFP should schedule firmware update according to package instructions"
$ws.Range("C26").Value = "MFP with SCC enabled"
$ws.Range("D26").Value = "Verify Install Package for policy violation"
$ws.Range("E26").Value = "1. Configure server to send policy violation package`n2. Trigger regular communication`n3. Monitor installation process"
$ws.Range("F26").Value = "MFP should apply policy settings as specified in package"
$ws.Range("C27").Value = "MFP with SCC enabled"
$ws.Range("D27").Value = "Verify Install Package for clone data"
$ws.Range("E27").Value = "1. Configure server to send clone data package`n2. Trigger regular communication`n3. Monitor installation process"
$ws.Range("F27").Value = "MFP should apply clone data as specified in package"
$ws.Range("C28").Value = "MFP with SCC enabled"
$ws.Range("D28").Value = "Verify Update Status functionality"
$ws.Range("E28").Value = "1. Configure server to send update package`n2. Trigger regular communication`n3. Monitor update status reporting"
$ws.Range("F28").Value = "MFP should report correct update status to server"
$ws.Range("C29").Value = "MFP with SCC enabled"
$ws.Range("D29").Value = "Verify Send Regular Data functionality"
$ws.Range("E29").Value = "1. Trigger regular communication`n2. Monitor regular data transmission"
$ws.Range("F29").Value = "MFP should collect and send all required data according to SendDataConfig setting"
$ws.Range("C30").Value = "MFP with SCC enabled"
$ws.Range("D30").Value = "Verify panel message display during SCC processing"
$ws.Range("F30").Value = "Panel should display `"Service in progress. Please do not turn off: XX`" with appropriate status code"
$ws.Range("C31").Value = "MFP with SCC enabled"
$ws.Range("D31").Value = "Verify event notification for device errors"
$ws.Range("E31").Value = "1. Generate device error with error code`n2. Monitor error notification"
$ws.Range("F31").Value = "MFP should send error notification to server"
$ws.Range("C32").Value = "MFP with SCC enabled"
$ws.Range("D32").Value = "Verify error resolution notification"
$ws.Range("E32").Value = "1. Generate device error`n2. Resolve error without power off`n3. Monitor error notification"
$ws.Range("F32").Value = "MFP should send error resolution notification with `"-`" prefix (e.g., -D102)"
$ws.Range("C33").Value = "MFP with SCC enabled"
$ws.Range("D33").Value = "Verify duplicate error handling"
$ws.Range("E33").Value = "1. Generate same device error multiple times`n2. Monitor error notifications"
$ws.Range("F33").Value = "MFP should not send duplicate error notifications unless error is resolved"
$ws.Range("C34").Value = "MFP with SCC enabled"
$ws.Range("D34").Value = "Verify Super Sleep interaction"
$ws.Range("E34").Value = "1. Put MFP in Super Sleep state`n2. Wait for scheduled communication time`n3. Monitor MFP behavior"
$ws.Range("F34").Value = "MFP should wake from Super Sleep, perform communication, then return to Super Sleep"
$ws.Range("C35").Value = "MFP with SCC enabled"
$ws.Range("D35").Value = "Verify Hibernation interaction"
$ws.Range("E35").Value = "1. Trigger SCC communication`n2. Attempt to put MFP in hibernation`n3. Monitor MFP behavior"
$ws.Range("F35").Value = "MFP should not enter hibernation during communication cycle"
$ws.Range("H35").Value = "Exception: during scheduled firmware update state"
$ws.Range("C36").Value = "MFP with SCC enabled"
$ws.Range("D36").Value = "Verify Persistent Policy functionality"
$ws.Range("E36").Value = "1. Configure server to set persistent policy`n2. Change settings locally`n3. Wait for persistent check interval`n4. Verify settings"
$ws.Range("F36").Value = "Settings should be restored according to persistent policy"
$ws.Range("C37").Value = "MFP with SCC enabled"
$ws.Range("D37").Value = "Verify Time-based Device State Data Values"
$ws.Range("E37").Value = "1. Configure server to set time-based values`n2. Monitor settings at different times"
$ws.Range("F37").Value = "Settings should change according to time-based configuration"
$ws.Range("C38").Value = "MFP with SCC enabled"
$ws.Range("D38").Value = "Verify IP Redirect functionality"
$ws.Range("E38").Value = "1. Configure MFP for first-time connection`n2. Monitor GetRedirectURL request`n3. Verify MFP connects to redirected URL"
$ws.Range("F38").Value = "MFP should connect to appropriate regional server based on IP address"
$ws.Range("C39").Value = "MFP with SCC enabled"
$ws.Range("D39").Value = "Verify URL forward setting"
$ws.Range("E39").Value = "1. Set 08-3827 to valid URL`n2. Set 08-3828 to 0 (Disable)`n3. Trigger SCC communication`n4. Monitor connection behavior"
$ws.Range("F39").Value = "MFP should not redirect to other servers even if server returns redirect URL"
$ws.Range("C40").Value = "MFP with SCC enabled"
$ws.Range("D40").Value = "Verify retry processing for connection errors"
$ws.Range("E40").Value = "1. Simulate connection error`n2. Monitor retry behavior"
$ws.Range("F40").Value = "MFP should retry in 60 seconds, then enter retry mode if unsuccessful"
$ws.Range("C41").Value = "MFP with SCC enabled"
$ws.Range("D41").Value = "Verify retry mode behavior"
$ws.Range("E41").Value = "1. Force MFP into retry mode`n2. Monitor communication attempts"
$ws.Range("F41").Value = "MFP should attempt communication once per day until successful"
$ws.Range("C42").Value = "MFP with SCC enabled"
$ws.Range("D42").Value = "Verify invalid URL handling"
$ws.Range("E42").Value = "1. Set 08-3827 to invalid URL`n2. Trigger SCC communication`n3. Monitor error handling"
$ws.Range("F42").Value = "MFP should report network error and not access default server"
$ws.Range("C43").Value = "MFP with SCC enabled"
$ws.Range("D43").Value = "Verify SSL certificate handling"
$ws.Range("E43").Value = "1. Configure custom HTTPS URL without importing certificate`n2. Trigger SCC communication`n3. Monitor error handling"
$ws.Range("F43").Value = "Communication should fail due to missing certificate"
$ws.Range("C44").Value = "MFP with SCC enabled"
$ws.Range("D44").Value = "Verify HDD full handling during download"
$ws.Range("E44").Value = "1. Fill HDD to capacity`n2. Trigger package download`n3. Monitor error handling"
$ws.Range("C45").Value = "MFP with SCC enabled"
$ws.Range("D45").Value = "Verify power failure recovery during download"
$ws.Range("E45").Value = "1. Trigger package download`n2. Simulate power failure during download`n3. Power on MFP`n4. Monitor recovery behavior"
$ws.Range("F45").Value = "MFP should delete incomplete downloaded data before next communication"
$ws.Range("B46").Value = "TC041"
$ws.Range("C46").Value = "MFP with SCC enabled"
$ws.Range("D46").Value = "Verify unzip failure handling"
$ws.Range("E46").Value = "1. Configure server to send corrupted ZIP package`n2. Trigger download`n3. Monitor error handling"
$ws.Range("F46").Value = "MFP should delete unzipped data and exit communication cycle"
$ws.Range("G46").Value = ""
$ws.Range("B47").Value = "TC042"
$ws.Range("C47").Value = "MFP with SCC enabled"
$ws.Range("D47").Value = "Verify firmware update scheduling"
$ws.Range("E47").Value = "1. Configure server to send firmware update with future schedule`n2. Trigger download`n3. Monitor scheduled update"
$ws.Range("F47").Value = "Firmware update should occur at scheduled time"
$ws.Range("G47").Value = ""
$ws.Range("B48").Value = "TC043"
$ws.Range("C48").Value = "MFP with SCC enabled"
$ws.Range("D48").Value = "Verify firmware update failure handling"
$ws.Range("E48").Value = "1. Configure server to send incompatible firmware`n2. Trigger download`n3. Monitor error handling"
$ws.Range("F48").Value = "MFP should report failure to server and schedule next communication"
$ws.Range("G48").Value = ""
$ws.Range("H48").Value = "After 3 failures, server should stop sending update"
